# Simplify the example data: rename PersonID values from "00N" to "ID00N"
# on both the DATA sheet and the PATIENTS sheet.

$wb = $excel.ActiveWorkbook

# ---- DATA sheet: column A holds PersonID, repeated across multiple rows ----
$dataSheet = $wb.Worksheets.Item("DATA")

$dataPersonRows = @{
    2  = "ID001"
    3  = "ID001"
    4  = "ID001"
    5  = "ID001"
    6  = "ID002"
    7  = "ID002"
    8  = "ID002"
    9  = "ID002"
    10 = "ID002"
    11 = "ID003"
    12 = "ID003"
    13 = "ID003"
    14 = "ID004"
    15 = "ID004"
    16 = "ID005"
    17 = "ID005"
    18 = "ID006"
    19 = "ID006"
    20 = "ID006"
}

foreach ($row in $dataPersonRows.Keys) {
    $dataSheet.Range("A$row").Value = $dataPersonRows[$row]
}

# ---- PATIENTS sheet: column A holds the unique PersonID per patient ----
$patientsSheet = $wb.Worksheets.Item("PATIENTS")

$patientPersonRows = @{
    2 = "ID001"
    3 = "ID002"
    4 = "ID003"
    5 = "ID004"
    6 = "ID005"
    7 = "ID006"
    8 = "ID007"
    9 = "ID008"
}

foreach ($row in $patientPersonRows.Keys) {
    $patientsSheet.Range("A$row").Value = $patientPersonRows[$row]
}

# ---- Reproduce the resulting cursor/selection positions ----
$dataSheet.Activate()
$dataSheet.Range("C23").Select()

$patientsSheet.Activate()
$patientsSheet.Range("A10").Select()

$dataSheet.Activate()
